# Auto-generated Excel COM-interop script to apply the cryptos.xlsx data refresh
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 16
$ws.Range("B16").Value = "Chainlink"
$ws.Range("C16").Value = "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "19.46"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "  +1.04%  "

# Row 17
$ws.Range("B17").Value = "Polygon"
$ws.Range("C17").Value = "https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "1.15"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "  +3.67%  "

# Row 27
$ws.Range("B27").Value = "LEO"
$ws.Range("C27").Value = "https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "6.06"
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "  +1.35%  "

# Row 28
$ws.Range("B28").Value = "Toncoin"
$ws.Range("C28").Value = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "3.79"
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = "  +2.28%  "

# Row 32
$ws.Range("B32").Value = "Cosmos"
$ws.Range("C32").Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "12.73"
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = "  +2.87%  "

# Row 33
$ws.Range("B33").Value = "InjectiveProtocol"
$ws.Range("C33").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "46.62"
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = "  +9.08%  "

# Row 2
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "68.778.17"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "  +1.69%  "

# Row 3
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.713.42"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "  +1.25%  "

# Row 4
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.00"
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "  -0.22%  "

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "620.22"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "  +8.62%  "

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "194.62"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "  +13.31%  "

# Row 7
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.637"
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "  +2.57%  "

# Row 8
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "1.00"
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "  -0.47%  "

# Row 9
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.727"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "  +4.35%  "

# Row 10
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "  +0.16%  "

# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "59.94"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "  +17.99%  "

# Row 12
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.0000288"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "  -0.39%  "

# Row 13
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "10.46"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "  +0.70%  "

# Row 14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "4.310.74"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "  +0.98%  "

# Row 15
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "3.717.98"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "  +1.27%  "

# Row 18
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "  +1.01%  "

# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "12.89"
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "  +0.92%  "

# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "68.670.71"
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "  +1.48%  "

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "410.94"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "  +1.98%  "

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "4.69"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "  +5.69%  "

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "90.24"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "  +3.72%  "

# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "3.09"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "  +2.23%  "

# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "11.58"
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "  +9.90%  "

# Row 26
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "13.09"
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "  +3.36%  "

# Row 29
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "9.68"
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = "  +3.37%  "

# Row 30
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "32.84"
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = "  +1.40%  "

# Row 31
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "7.65"
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = "  +2.43%  "

# Row 34
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.123"
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = "  +7.03%  "

# Row 35
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "635.96"
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = "  +6.65%  "

# Row 36
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "67.48"
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = "  +4.31%  "

# Row 37
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = "  +5.73%  "

# Row 38
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.0₃0827"
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = "  -6.48%  "

# Row 39
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.999"
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "  -0.17%  "

# Row 40
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "  -0.09%  "

# Row 41
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.140"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "  +5.51%  "

# Row 42
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "3.04"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "  +3.13%  "

# Row 43
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.0448"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "  +3.31%  "

# Row 44
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "2.62"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "  +1.05%  "

# Row 45
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "2.941.95"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "  +7.62%  "

# Row 46
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "9.50"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "  +4.35%  "

# Row 47
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.139"
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "  +4.92%  "

# Row 48
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "  +0.37%  "

# Row 49
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "145.94"
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "  +1.97%  "

# Row 50
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "3.07"
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "  -1.03%  "

# Row 51
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "2.77"
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "  +2.02%  "

